$wb = $excel.ActiveWorkbook

$ws1 = $wb.Sheets("About")
$ws2 = $wb.Sheets("PDiBCpDoC")

# --- "About" sheet -------------------------------------------------------

# Remove the embedded chart image - source is no longer the BNEF chart.
if ($ws1.Shapes.Count -gt 0) {
    for ($i = $ws1.Shapes.Count; $i -ge 1; $i--) {
        $ws1.Shapes.Item($i).Delete()
    }
}

# Update the source citation to the new paper (MIT / RSC publication).
# Order chosen so new shared strings land in the same slot order as the
# target workbook: URL, then paper title, then publisher/author.
$ws2.Range("B2").Formula = "=AVERAGE(0.2,0.27)"

$ws1.Range("B6").Value = "https://pubs.rsc.org/en/content/articlepdf/2021/ee/d0ee02681f?page=search"
$ws1.Range("B5").Value = "Re-examining rates of lithium-ion battery technology improvement and cost decline"
$ws1.Range("B3").Value = "Massachusetts Institute of Technology"
$ws1.Range("B4").Value = 2021
$ws1.Range("B7").Value = "Abstract"

# The old note below the chart ("the graph only extends to 2030...") no
# longer applies - clear it but keep the italic style on C8.
$ws1.Range("C8").ClearContents()

# New footnote about how the learning rate was derived.
$ws1.Range("A9").Value = "Note: We take the average of learning rates quoted in the Abstract (20%-27%)"

# --- "PDiBCpDoC" sheet ----------------------------------------------------
# Perc Decline per Doubling is now the average of the 20%-27% learning-rate
# range quoted in the new source, instead of a hard-coded 0.18.
$ws2.Range("I4").Select() | Out-Null

# Re-activate "About" (keep it the selected/front-most tab) and park the
# cursor below the new footnote, matching where the author left off.
$ws1.Range("A10").Select() | Out-Null
